# Insert two new data rows right after row 172 (before the former row 173),
# which pushes the former rows 173..212 down to 175..214.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A173:R174").EntireRow.Insert()

# New row 173: Camote, "1a (cosecha)" entry dated 44641
$ws.Cells.Item(173, 1).Value  = 5
$ws.Cells.Item(173, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(173, 3).Value  = "Maule"
$ws.Cells.Item(173, 4).Value2 = 44641
$ws.Cells.Item(173, 5).Value  = 7
$ws.Cells.Item(173, 6).Value  = 100112045
$ws.Cells.Item(173, 7).Value  = "Zapallo"
$ws.Cells.Item(173, 8).Value  = "Camote"
$ws.Cells.Item(173, 9).Value  = "1a (cosecha)"
$ws.Cells.Item(173, 10).Value = 900
$ws.Cells.Item(173, 11).Value = 300
$ws.Cells.Item(173, 12).Value = 300
$ws.Cells.Item(173, 13).Value = 300
$ws.Cells.Item(173, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(173, 15).Value = "Región del Maule"
$ws.Cells.Item(173, 16).Value = 300
$ws.Cells.Item(173, 17).Value = 1
$ws.Cells.Item(173, 18).Value = "Hortaliza"

# New row 174: Paine, "1a (guarda)" entry dated 44641
$ws.Cells.Item(174, 1).Value  = 5
$ws.Cells.Item(174, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(174, 3).Value  = "Maule"
$ws.Cells.Item(174, 4).Value2 = 44641
$ws.Cells.Item(174, 5).Value  = 7
$ws.Cells.Item(174, 6).Value  = 100112045
$ws.Cells.Item(174, 7).Value  = "Zapallo"
$ws.Cells.Item(174, 8).Value  = "Paine"
$ws.Cells.Item(174, 9).Value  = "1a (guarda)"
$ws.Cells.Item(174, 10).Value = 2000
$ws.Cells.Item(174, 11).Value = 100
$ws.Cells.Item(174, 12).Value = 100
$ws.Cells.Item(174, 13).Value = 100
$ws.Cells.Item(174, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(174, 15).Value = "Región del Maule"
$ws.Cells.Item(174, 16).Value = 100
$ws.Cells.Item(174, 17).Value = 1
$ws.Cells.Item(174, 18).Value = "Hortaliza"
